$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" (E), "Valor Mora" (F) and "Salario Basico" (G) for
# the data rows 16-20. The account-statement periods now run in ascending
# order (2405..2409) instead of descending (2409..2405), and the Mora /
# Salario values for the new periods are refreshed.
$ws.Range("E16").Value = "2405"
$ws.Range("F16").Value = 56144
$ws.Range("G16").Value = 1754700

$ws.Range("E17").Value = "2406"
$ws.Range("F17").Value = 70180
$ws.Range("G17").Value = 1754700

$ws.Range("E18").Value = "2407"
$ws.Range("F18").Value = 70180
$ws.Range("G18").Value = 1754700

$ws.Range("E19").Value = "2408"
$ws.Range("F19").Value = 70180
$ws.Range("G19").Value = 1754700

$ws.Range("E20").Value = "2409"
$ws.Range("F20").Value = 70180
$ws.Range("G20").Value = 1754700
